# This script reproduces the target edit:
#  - Two new data rows are inserted into the worksheet right before the
#    existing row 569 (shifting the old rows 569..677 down to 571..679).
#  - The two newly inserted rows (569 and 570) contain new price-report
#    entries, built by copying the former row 569 / row 570 entries and
#    changing a handful of fields (date, volume, price range, origin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 569. Everything that used to
# be on rows 569-677 shifts down to rows 571-679.
$ws.Rows("569:570").Insert()

# --- New row 569 -----------------------------------------------------
# Copy of the (now shifted) row 571 content, with Fecha and Volumen updated.
$ws.Range("A569").Value = 10
$ws.Range("B569").Value = "Vega Modelo de Temuco"
$ws.Range("C569").Value = "La Araucanía"
$ws.Range("D569").Value = 45209
$ws.Range("E569").Value = 9
$ws.Range("F569").Value = 100114014
$ws.Range("G569").Value = "Betarraga"
$ws.Range("H569").Value = "Sin especificar"
$ws.Range("I569").Value = "Primera"
$ws.Range("J569").Value = 70
$ws.Range("K569").Value = 8000
$ws.Range("L569").Value = 8000
$ws.Range("M569").Value = 8000
$ws.Range("N569").Value = "$/docena de paquetes"
$ws.Range("O569").Value = "Provincia de Cautín"
$ws.Range("P569").Value = 667
$ws.Range("Q569").Value = 12
$ws.Range("R569").Value = "Hortaliza"

# --- New row 570 -----------------------------------------------------
# Copy of the (now shifted) row 572 content, with Fecha, Volumen, Precio
# máximo/promedio, Origen and Precio $/Kg updated.
$ws.Range("A570").Value = 10
$ws.Range("B570").Value = "Vega Modelo de Temuco"
$ws.Range("C570").Value = "La Araucanía"
$ws.Range("D570").Value = 45209
$ws.Range("E570").Value = 9
$ws.Range("F570").Value = 100114014
$ws.Range("G570").Value = "Betarraga"
$ws.Range("H570").Value = "Sin especificar"
$ws.Range("I570").Value = "Primera"
$ws.Range("J570").Value = 80
$ws.Range("K570").Value = 8000
$ws.Range("L570").Value = 8000
$ws.Range("M570").Value = 8000
$ws.Range("N570").Value = "$/docena de paquetes"
$ws.Range("O570").Value = "Región del Maule"
$ws.Range("P570").Value = 667
$ws.Range("Q570").Value = 12
$ws.Range("R570").Value = "Hortaliza"
